# A new price-record row is inserted at row 554 (new entry for Betarraga,
# Macroferia Regional de Talca). All existing rows from 554 downward shift
# down by one (554->555, ..., 671->672), and the new row 554 is populated
# with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 554, pushing rows 554..671 down to 555..672
$ws.Rows.Item(554).Insert()

# Populate the new row 554 with the inserted record's values
$ws.Cells.Item(554, 1).Value = 5
$ws.Cells.Item(554, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(554, 3).Value = "Maule"
$ws.Cells.Item(554, 4).Value = Get-Date -Year 2023 -Month 10 -Day 12 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(554, 5).Value = 7
$ws.Cells.Item(554, 6).Value = 100114014
$ws.Cells.Item(554, 7).Value = "Betarraga"
$ws.Cells.Item(554, 8).Value = "Sin especificar"
$ws.Cells.Item(554, 9).Value = "Primera"
$ws.Cells.Item(554, 10).Value = 3000
$ws.Cells.Item(554, 11).Value = 500
$ws.Cells.Item(554, 12).Value = 500
$ws.Cells.Item(554, 13).Value = 500
$ws.Cells.Item(554, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(554, 15).Value = "Región del Maule"
$ws.Cells.Item(554, 16).Value = 100
$ws.Cells.Item(554, 17).Value = 5
$ws.Cells.Item(554, 18).Value = "Hortaliza"
